$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (No.1 - Abdul Latief Baedhowi): remove TTL place/birthdate info and reset UMUR
$ws.Range("B4").Value = "Abdul Latief Baedhowi//"
$ws.Range("F4").Value = 0

# Row 5 (No.2 - Ahmad Hanif): clear Jabatan (Admin removed) and change Pendidikan S1 -> S2
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = "S2"

# Row 6 (No.3 - Iwan Setiyawan): add TTL info, set Pendidikan to S1
$ws.Range("B6").Value = "Iwan Setiyawan/Bantul/2022-06-06"
$ws.Range("E6").Value = "S1"

# Rows with D3 -> DIII rename (Roni Slamet, Ihsanuddin, Bambang Gunartok)
$ws.Range("E12").Value = "DIII"
$ws.Range("E13").Value = "DIII"
$ws.Range("E15").Value = "DIII"

# Row 20 (Bagus Widodo): set Pendidikan to DIII
$ws.Range("E20").Value = "DIII"

# Row 31 (Nur Muh Attabik): set Pendidikan to S1
$ws.Range("E31").Value = "S1"

# Row 148 (Ihsan Firdaus): set Pendidikan to DIII
$ws.Range("E148").Value = "DIII"

# Row 168 (Burhanu Sultan Ramadan): set Jabatan to Staff
$ws.Range("D168").Value = "Staff"
